# Updated data import: each cross row now records an explicit "F" (sex)
# value in column D, and the stray scratch-note cells (I2, I3, K7) that
# held old to-do comments are cleared out now that the file points at
# the new working/rego source data instead of data/.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D7").Value = "F"

$ws.Range("I2").Value = $null
$ws.Range("I3").Value = $null
$ws.Range("K7").Value = $null

# Restore the cursor/selection to match the refreshed view.
[void]$ws.Range("H1:M8").Select()
